# Update the cryptocurrency price/volume table on Sheet1 with the latest
# scraped values (GitHub Actions refresh). Rows 2-47 keep the same coin but
# get refreshed Price (D) / Volume(1h) (E) figures; rows 48-51 shift down by
# one because a new coin (BabyDogeCoin) was inserted at row 48, pushing the
# previous USDD row off the bottom of the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.869.58'
$ws.Range('E2').Value = '  +0.61%  '

$ws.Range('D3').Value = '1.641.55'
$ws.Range('E3').Value = '  +0.06%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '218.43'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.05%  '

$ws.Range('E6').Value = '  -0.73%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = "Normal"

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.251'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.13%  '

$ws.Range('E9').Value = '  -0.59%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.27'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.99%  '

$ws.Range('E11').Value = '  +0.21%  '

$ws.Range('D12').Value = '1.869.98'
$ws.Range('E12').Value = '  +0.09%  '

$ws.Range('D13').Value = '1.647.58'
$ws.Range('E13').Value = '  +1.01%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.15'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.37%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.528'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.35%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '65.31'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +1.39%  '

$ws.Range('D17').Value = '26.861.01'
$ws.Range('E17').Value = '  +0.54%  '

$ws.Range('E18').Value = '  -0.37%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '215.68'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.60%  '

$ws.Range('E20').Value = '  -0.07%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.37'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.05%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.54'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +4.74%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.38'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -2.35%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '9.21'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -1.21%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '147.46'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.25%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.64%  '

$ws.Range('E27').Value = '  -0.34%  '

$ws.Range('E28').Value = '  +1.29%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.75'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.83%  '

$ws.Range('E30').Value = '  -0.07%  '

$ws.Range('E31').Value = '  +0.97%  '

$ws.Range('E32').Value = '  +0.64%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.00'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.02%  '

$ws.Range('E34').Value = '  +1.73%  '

$ws.Range('D35').Value = '1.273.82'
$ws.Range('E35').Value = '  -1.22%  '

$ws.Range('E36').Value = '  +0.23%  '

$ws.Range('E37').Value = '  -1.79%  '

$ws.Range('E38').Value = '  -0.78%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.820'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.46%  '

$ws.Range('E40').Value = '  -0.02%  '

$ws.Range('E41').Value = '  +0.04%  '

$ws.Range('E42').Value = '  +0.03%  '

$ws.Range('E43').Value = '  -0.64%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '92.55'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.27%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '61.04'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.63%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.04'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -8.56%  '

$ws.Range('E47').Value = '  -0.05%  '

$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0102'
$ws.Range('E48').Value = '  -1.60%  '

$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0515'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -2.08%  '

$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0970'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.20%  '

$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.57'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.74%  '
